$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the new, longer article name
$ws.Columns(1).ColumnWidth = 25.666666666666668

# Add the new "Elvis Presley" row of data (row 3)
$ws.Range("A3").Value = "Elvis Presley (Largest in English src: https://diff.wikimedia.org/2016/05/12/rock-n-scroll-english-wikipedias-longest-featured-articles/)"
$ws.Range("B3").Value = "en"
$ws.Range("D3").Value = 589

# Match the author's final cursor position
$ws.Range("B10").Select()
